$d = $word.ActiveDocument

# Locate the first paragraph ("Design Document Outline"). We will insert three
# new paragraphs immediately after it: an empty paragraph, a paragraph holding
# an italicized pull-quote, and another empty paragraph.
$p1 = $d.Paragraphs.Item(1)

# Build a range that ends just before paragraph 1's trailing paragraph mark so
# that inserting XML there appends new paragraphs after it instead of
# replacing/merging into the existing paragraph mark.
$insertionPoint = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$insertionPoint.Collapse(0)

$quoteText = '"The world outside the dungeon isn’t designed to be a caricature—it doesn’t need exaggeration to reveal its absurdity, horror, and humor. NPCs aren’t shallow tropes, but realistic figures with motivations that make them feel alive—whether as villains, opportunists, or simply desperate souls making their way. The contrast is clear: the dungeon is danger, but predictable; the outside world is survival, but increasingly suffocating. There’s always a choice: risk the wild unknown or endure the exhausting routine. Neither is fair, but one feels like freedom."'

$bodyFragment = '<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:bidi w:val="0"/><w:spacing w:lineRule="auto" w:line="276" w:before="0" w:after="140"/><w:jc w:val="left"/><w:rPr><w:rStyle w:val="Strong"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/><w:bidi w:val="0"/><w:spacing w:lineRule="auto" w:line="276" w:before="0" w:after="140"/><w:jc w:val="left"/><w:rPr><w:rStyle w:val="Strong"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="Strong"/><w:i/><w:iCs/></w:rPr><w:t>' + $quoteText + '</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/><w:bidi w:val="0"/><w:spacing w:lineRule="auto" w:line="276" w:before="0" w:after="140"/><w:jc w:val="left"/><w:rPr><w:rStyle w:val="Strong"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p>'

$flatOpc = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>__BODY__</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$flatOpc = $flatOpc.Replace('__BODY__', $bodyFragment)

$insertionPoint.InsertXML($flatOpc)
